$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 16.31496124559424
$ws.Range("C2").Value = 44.300831108300457
$ws.Range("D2").Value = 47.442058184454872
$ws.Range("E2").Value = 48.662927375082532

# Row 3 data values
$ws.Range("B3").Value = 27.130757766161199
$ws.Range("C3").Value = 68.091621566852538
$ws.Range("D3").Value = 57.162003373349634
$ws.Range("E3").Value = 45.27787645089181

# Update selection to reflect the edited range
$ws.Range("B1:E3").Select()
